$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D42").Value = "다전공제도, 다전공뭐야, 다전공무엇, 다전공알려줘, 다전공설명, 다전공개념, 다전공정의, 다전공이뭐"
$ws.Range("D43").Value = "복수전공복수전공, 복수전공제도, 복수전공뭐야, 복수전공무엇, 복수전공알려줘, 복수전공설명, 복수전공개념, 복수전공정의, 복수전공이뭐"
$ws.Range("D44").Value = "부전공부전공, 부전공제도, 부전공뭐야, 부전공무엇, 부전공알려줘, 부전공설명, 부전공개념, 부전공정의, 부전공이뭐"
$ws.Range("D45").Value = "융합전공융합전공, 융합전공제도, 융합전공뭐야, 융합전공무엇, 융합전공알려줘, 융합전공설명, 융합전공개념, 융합전공정의, 융합전공이뭐"
$ws.Range("D46").Value = "융합부전공융합부전공, 융합부전공제도, 융합부전공뭐야, 융합부전공무엇, 융합부전공알려줘, 융합부전공설명, 융합부전공개념, 융합부전공정의, 융합부전공이뭐"
$ws.Range("D47").Value = "연계전공연계전공, 연계전공제도, 연계전공뭐야, 연계전공무엇, 연계전공알려줘, 연계전공설명, 연계전공개념, 연계전공정의, 연계전공이뭐"
$ws.Range("D48").Value = "마이크로디그리란, 마이크로디그리가뭐, 마이크로디그리가뭐야, 마이크로디그리제도, 마이크로디그리설명, 마디란, 마디가뭐, 마이크로디그리무엇, 마이크로디그리알려, 마이크로디그리개념, 마이크로디그리정의, microdegree뭐, 소단위전공과정이란, 소단위전공과정뭐, 소단위전공과정이뭐, 소단위가뭐"

$wb.Save()
